$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format first so numeric-looking strings
# (e.g. "1.000", "317.16") are stored as text, matching the original
# inline-string cell type instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.048.75'
$ws.Range("E2").Value = '  -1.27%  '
$ws.Range("D3").Value = '1.793.21'
$ws.Range("E3").Value = '  -0.19%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '317.16'
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").Value = '0.5362'
$ws.Range("E7").Value = '  -1.70%  '
$ws.Range("D8").Value = '0.3771'
$ws.Range("E8").Value = '  -1.62%  '
$ws.Range("D9").Value = '0.07428'
$ws.Range("E9").Value = '  -2.27%  '
$ws.Range("D10").Value = '41.80'
$ws.Range("E10").Value = '  -1.97%  '
$ws.Range("D11").Value = '1.092'
$ws.Range("E11").Value = '  -2.88%  '
$ws.Range("D12").Value = '1.000'
$ws.Range("E12").Value = '  -0.14%  '
$ws.Range("D13").Value = '20.56'
$ws.Range("E13").Value = '  -2.78%  '
$ws.Range("D14").Value = '6.122'
$ws.Range("E14").Value = '  -1.08%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '7.243'
$ws.Range("E15").Value = '  -2.20%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '1.785.38'
$ws.Range("E16").Value = '  -0.53%  '
$ws.Range("D17").Value = '88.99'
$ws.Range("E17").Value = '  -2.79%  '
$ws.Range("E18").Value = '  -1.25%  '
$ws.Range("D19").Value = '0.06502'
$ws.Range("E19").Value = '  +0.85%  '
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  +0.04%  '
$ws.Range("D21").Value = '17.26'
$ws.Range("E21").Value = '  -0.32%  '
$ws.Range("D22").Value = '5.895'
$ws.Range("E22").Value = '  -1.13%  '
$ws.Range("D23").Value = '28.062.48'
$ws.Range("E23").Value = '  -1.24%  '
$ws.Range("E24").Value = '  -2.71%  '
$ws.Range("E25").Value = '  -1.91%  '
$ws.Range("D26").Value = '155.47'
$ws.Range("E26").Value = '  -2.64%  '
$ws.Range("D27").Value = '20.29'
$ws.Range("E27").Value = '  -2.03%  '
$ws.Range("D28").Value = '1.992.84'
$ws.Range("E28").Value = '  -0.44%  '
$ws.Range("D29").Value = '2.301'
$ws.Range("E29").Value = '  -3.99%  '
$ws.Range("D30").Value = '121.22'
$ws.Range("E30").Value = '  -1.76%  '
$ws.Range("D31").Value = '1.117'
$ws.Range("E31").Value = '  -0.86%  '
$ws.Range("D32").Value = '0.1059'
$ws.Range("E32").Value = '  +3.01%  '
$ws.Range("D33").Value = '3.663'
$ws.Range("E33").Value = '  -0.57%  '
$ws.Range("D34").Value = '5.558'
$ws.Range("E34").Value = '  -3.54%  '
$ws.Range("B35").Value = 'Algorand'
$ws.Range("C35").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D35").Value = '0.2248'
$ws.Range("E35").Value = '  -5.09%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = '0.06489'
$ws.Range("E36").Value = '  -4.64%  '
$ws.Range("E37").Value = '  -1.46%  '
$ws.Range("D38").Value = '5.017'
$ws.Range("E38").Value = '  -2.91%  '
$ws.Range("D39").Value = '8.461'
$ws.Range("E39").Value = '  -3.50%  '
$ws.Range("D40").Value = '0.6193'
$ws.Range("E40").Value = '  -3.21%  '
$ws.Range("D41").Value = '1.450'
$ws.Range("E41").Value = '  +2.82%  '
$ws.Range("D42").Value = '11.14'
$ws.Range("E42").Value = '  -4.61%  '
$ws.Range("D43").Value = '1.177'
$ws.Range("E43").Value = '  +1.52%  '
$ws.Range("D44").Value = '0.9998'
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("D45").Value = '13.24'
$ws.Range("E45").Value = '  -2.52%  '
$ws.Range("E46").Value = '  +0.20%  '
$ws.Range("D47").Value = '0.5780'
$ws.Range("E47").Value = '  -3.31%  '
$ws.Range("D48").Value = '125.05'
$ws.Range("E48").Value = '  -1.18%  '
$ws.Range("D49").Value = '1.187'
$ws.Range("E49").Value = '  +3.18%  '
$ws.Range("D50").Value = '1.926'
$ws.Range("E50").Value = '  -3.72%  '
$ws.Range("D51").Value = '0.06817'
$ws.Range("E51").Value = '  -1.95%  '

# Restore the default "Normal" style on column D so no residual
# number-format styling is left on the cells (matches original OOXML).
$ws.Range("D2:D51").Style = "Normal"
